$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "14×63=882" "17×51=867"
Replace-Text "24×61=1464" "74×49=3626"
Replace-Text "77×38=2926" "68×25=1700"
Replace-Text "65×47=3055" "55×45=2475"
Replace-Text "72×13=936" "63×56=3528"
Replace-Text "23×95=2185" "31×56=1736"
Replace-Text "71×53=3763" "13×48=624"
Replace-Text "31×21=651" "88×79=6952"
Replace-Text "97×28=2716" "98×17=1666"
Replace-Text "82×51=4182" "55×40=2200"
Replace-Text "17×54=918" "47×83=3901"
Replace-Text "24×95=2280" "75×56=4200"
Replace-Text "98×92=9016" "92×60=5520"
Replace-Text "70×71=4970" "44×45=1980"
Replace-Text "96×82=7872" "97×25=2425"
Replace-Text "78×26=2028" "68×51=3468"
Replace-Text "95×84=7980" "45×49=2205"
Replace-Text "85×68=5780" "37×98=3626"
Replace-Text "68×76=5168" "44×79=3476"
Replace-Text "93×49=4557" "89×39=3471"
Replace-Text "78×51=3978" "70×73=5110"
Replace-Text "27×30=810" "75×24=1800"
Replace-Text "13×56=728" "39×44=1716"
Replace-Text "11×42=462" "21×12=252"
Replace-Text "14×93=1302" "99×76=7524"
